$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '35.351.41'
$ws.Range('E2').Value = '  +0.44%  '

# Row 3
$ws.Range('D3').Value = '1.913.28'
$ws.Range('E3').Value = '  +2.91%  '

# Row 4
$ws.Range('E4').Value = '  -0.49%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.99%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.657'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.59%  '

# Row 7
$ws.Range('E7').Value = '  -0.46%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.26'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.26%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.350'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.49%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.77'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.49%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0715'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.10%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0995'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.51%  '

# Row 13
$ws.Range('D13').Value = '2.194.58'
$ws.Range('E13').Value = '  +3.05%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.08%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.700'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.64%  '

# Row 16
$ws.Range('D16').Value = '1.917.26'
$ws.Range('E16').Value = '  +3.05%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.28%  '

# Row 18
$ws.Range('D18').Value = '35.334.18'
$ws.Range('E18').Value = '  +0.46%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.21%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0829'
$ws.Range('E20').Value = '  +4.17%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '239.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.40%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.29%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.97%  '

# Row 24
$ws.Range('E24').Value = '  -0.49%  '

# Row 25
$ws.Range('E25').Value = '  +1.09%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +21.25%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.92%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.54%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.44%  '

# Row 30
$ws.Range('E30').Value = '  +2.24%  '

# Row 31
$ws.Range('E31').Value = '  +3.80%  '

# Row 32
$ws.Range('E32').Value = '  +1.94%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.937'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.72%  '

# Row 34
$ws.Range('E34').Value = '  -0.39%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.62%  '

# Row 36
$ws.Range('E36').Value = '  -3.56%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.00%  '

# Row 38
$ws.Range('E38').Value = '  +0.98%  '

# Row 39
$ws.Range('E39').Value = '  +2.26%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0661'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +13.55%  '

# Row 41
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.65%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0208'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.96%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '90.32'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.12%  '

# Row 44
$ws.Range('D44').Value = '1.341.78'
$ws.Range('E44').Value = '  +0.13%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.55%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '47.98'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +38.07%  '

# Row 47
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.14%  '

# Row 48
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.75%  '

# Row 49
$ws.Range('E49').Value = '  +0.48%  '

# Row 50
$ws.Range('D50').Value = '2.100.37'
$ws.Range('E50').Value = '  +2.80%  '

# Row 51
$ws.Range('E51').Value = '  +3.10%  '
